# Append the 09/29/2025 profit-allocation row (A28:C28) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A holds the date as literal text (matching every other row in the
# sheet, which stores "MM/DD/YYYY" as a string rather than a real date).
# Force text storage so Excel doesn't auto-convert the date-looking string
# into a date serial number, then drop the temporary number format so the
# cell ends up unstyled like its neighbours.
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = "09/29/2025"
$ws.Range("A28").ClearFormats()

$ws.Range("B28").Value = 0.1388883331783817
$ws.Range("C28").Value = 0.8611116668216183
